$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet SCD0193 -> SCD0011
$ws.Name = "SCD0011"

# Update TC_ID cell (B2): DGS-208 -> SCD0011-024
$ws.Range("B2").Value = "SCD0011-024"

# Align the used range (A1:Z2) to the left horizontally (and keep/apply
# vertical centering so previously-unstyled cells match the rest of the row)
$ws.Range("A1:Z2").HorizontalAlignment = -4131
$ws.Range("A1:Z2").VerticalAlignment = -4108

# Scroll the sheet view back so the top-left visible cell is A1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
